$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.710.11"
$ws.Range("E2").Value = "  +6.91%  "
$ws.Range("D3").Value = "1.810.06"
$ws.Range("E3").Value = "  +4.75%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06395"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("D10").Value = "1.804.76"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07108"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6503"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.719"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "81.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").Value = "28.676.73"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9992"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007362"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9993"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.26%  "
$ws.Range("D21").Value = "2.034.81"
$ws.Range("E21").Value = "  +4.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.620"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.906"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.319"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.888"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.392"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.193"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08375"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.844"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04975"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.093"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6746"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.664"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9625"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.662"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01593"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.930"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +6.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.237"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1226"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05494"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.180"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.11%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3608"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.76%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.302"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.03%  "
